$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column B (shifts old B..I to C..J)
$ws.Columns.Item(2).Insert()

# 2. Insert 3 new rows before row 3 (shifts old rows 3.. down by 3)
$ws.Rows.Item(3).Resize(3).Insert()

# 3. Clear the stray cells left behind in rows 6-8 column A by the row-insert shift
#    (these held the old "Incep variable / UNet / Baseline" single-row table which is
#    being replaced by the new Baseline/UNet/InceptionNet x Fixed/Variable block)
$ws.Range("A6").Value = ""
$ws.Range("A7").Value = ""
$ws.Range("A8").Value = ""
$ws.Range("B7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("F7").Value = ""

# 4. Populate the new "Model" sub-table (rows 2-9), typing new unique strings in the
#    same order they were first entered so shared-string ids line up.
$ws.Range("A7").Value = "InceptionNet"
$ws.Range("B2").Value = "Overlaps"
$ws.Range("B3").Value = "Fixed"
$ws.Range("B4").Value = "Variable"
$ws.Range("A2").Value = "Model"

$ws.Range("A3").Value = "Baseline"
$ws.Range("B5").Value = "Fixed"
$ws.Range("A5").Value = "UNet"
$ws.Range("B6").Value = "Variable"
$ws.Range("B7").Value = "Fixed"
$ws.Range("B8").Value = "Variable"

# 5. New numeric "running time" values
$ws.Range("I5").Value = 2.239
$ws.Range("I6").Value = 1.187
$ws.Range("I8").Value = 1.239

# 6. Centre-align column A for the new block (gives it style index 1, matching header style)
$ws.Range("A3:A8").HorizontalAlignment = -4108

# 7. Merge the vertically-stacked "Model" labels
$ws.Range("A3:A4").Merge()
$ws.Range("A5:A6").Merge()
$ws.Range("A7:A8").Merge()

# 8. Column widths (best-effort; engine cannot reproduce bestFit auto-measurement exactly)
$ws.Columns.Item(1).ColumnWidth = 12.6
$ws.Columns.Item(2).ColumnWidth = 12.6

# 9. Update the selection to match the saved workbook state
$ws.Range("I4").Select()
